$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.055.39'
$ws.Range('E2').Value = '  +2.42%  '
$ws.Range('D3').Value = '1.675.60'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0618'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.21'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.91%  '
$ws.Range('D12').Value = '1.912.36'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('D13').Value = '1.678.21'
$ws.Range('E13').Value = '  +3.64%  '
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.91'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.05%  '
$ws.Range('D17').Value = '27.086.17'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '237.09'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +3.98%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.69'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.16'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0499'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('D33').Value = '1.477.47'
$ws.Range('E33').Value = '  -3.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.33%  '
$ws.Range('E35').Value = '  +5.99%  '
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.576'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.901'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.91%  '
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('E41').Value = '  +12.22%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.28'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.70'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.05%  '
$ws.Range('D45').Value = '1.823.40'
$ws.Range('E45').Value = '  +3.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.780'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.37'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.69'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.87%  '
